$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 165
$ws1.Range("F6").Value = 9419
$ws1.Range("F9").Value = 1204
$ws1.Range("F10").Value = 1148
$ws1.Range("F14").Value = 262
$ws1.Range("F15").Value = 426
$ws1.Range("F16").Value = 89
$ws1.Range("F18").Value = 1286

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 165
$ws4.Range("F7").Value = 9419
$ws4.Range("F10").Value = 1204
$ws4.Range("F11").Value = 1148
$ws4.Range("F15").Value = 262
$ws4.Range("F16").Value = 426
$ws4.Range("F17").Value = 89
$ws4.Range("F19").Value = 1286
